$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 2023-09-05 (45174)
# to 2023-09-06 (45175) for rows 2 through 8.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
